$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Factorial design table: update Max_depth (F), Learning_Rate (G),
# N_estimators (H) and colsample_bytree (I) per row with the new
# Design Expert factor levels.
$data = @(
  @(2, 3, 0.010000000000000002, 100, 0.4),
  @(3, 3, 0.010000000000000002, 100, 0.4),
  @(4, 3, 0.010000000000000002, 100, 0.4),
  @(5, 3, 0.010000000000000002, 100, 0.4),
  @(6, 9, 0.010000000000000002, 100, 0.4),
  @(7, 9, 0.010000000000000002, 100, 0.4),
  @(8, 9, 0.010000000000000002, 100, 0.4),
  @(9, 9, 0.010000000000000002, 100, 0.4),
  @(10, 3, 0.09, 100, 0.4),
  @(11, 3, 0.09, 100, 0.4),
  @(12, 3, 0.09, 100, 0.4),
  @(13, 3, 0.09, 100, 0.4),
  @(14, 9, 0.09, 100, 0.4),
  @(15, 9, 0.09, 100, 0.4),
  @(16, 9, 0.09, 100, 0.4),
  @(17, 9, 0.09, 100, 0.4),
  @(18, 3, 0.010000000000000002, 800, 0.4),
  @(19, 3, 0.010000000000000002, 800, 0.4),
  @(20, 3, 0.010000000000000002, 800, 0.4),
  @(21, 3, 0.010000000000000002, 800, 0.4),
  @(22, 9, 0.010000000000000002, 800, 0.4),
  @(23, 9, 0.010000000000000002, 800, 0.4),
  @(24, 9, 0.010000000000000002, 800, 0.4),
  @(25, 9, 0.010000000000000002, 800, 0.4),
  @(26, 3, 0.09, 800, 0.4),
  @(27, 3, 0.09, 800, 0.4),
  @(28, 3, 0.09, 800, 0.4),
  @(29, 3, 0.09, 800, 0.4),
  @(30, 9, 0.09, 800, 0.4),
  @(31, 9, 0.09, 800, 0.4),
  @(32, 9, 0.09, 800, 0.4),
  @(33, 9, 0.09, 800, 0.4),
  @(34, 3, 0.010000000000000002, 100, 0.9),
  @(35, 3, 0.010000000000000002, 100, 0.9),
  @(36, 3, 0.010000000000000002, 100, 0.9),
  @(37, 3, 0.010000000000000002, 100, 0.9),
  @(38, 9, 0.010000000000000002, 100, 0.9),
  @(39, 9, 0.010000000000000002, 100, 0.9),
  @(40, 9, 0.010000000000000002, 100, 0.9),
  @(41, 9, 0.010000000000000002, 100, 0.9),
  @(42, 3, 0.09, 100, 0.9),
  @(43, 3, 0.09, 100, 0.9),
  @(44, 3, 0.09, 100, 0.9),
  @(45, 3, 0.09, 100, 0.9),
  @(46, 9, 0.09, 100, 0.9),
  @(47, 9, 0.09, 100, 0.9),
  @(48, 9, 0.09, 100, 0.9),
  @(49, 9, 0.09, 100, 0.9),
  @(50, 3, 0.010000000000000002, 800, 0.9),
  @(51, 3, 0.010000000000000002, 800, 0.9),
  @(52, 3, 0.010000000000000002, 800, 0.9),
  @(53, 3, 0.010000000000000002, 800, 0.9),
  @(54, 9, 0.010000000000000002, 800, 0.9),
  @(55, 9, 0.010000000000000002, 800, 0.9),
  @(56, 9, 0.010000000000000002, 800, 0.9),
  @(57, 9, 0.010000000000000002, 800, 0.9),
  @(58, 3, 0.09, 800, 0.9),
  @(59, 3, 0.09, 800, 0.9),
  @(60, 3, 0.09, 800, 0.9),
  @(61, 3, 0.09, 800, 0.9),
  @(62, 9, 0.09, 800, 0.9),
  @(63, 9, 0.09, 800, 0.9),
  @(64, 9, 0.09, 800, 0.9),
  @(65, 9, 0.09, 800, 0.9),
  @(66, 6, 0.05, 450, 0.65),
  @(67, 6, 0.05, 450, 0.65),
  @(68, 6, 0.05, 450, 0.65),
  @(69, 6, 0.05, 450, 0.65)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 6).Value = $row[1]
  $ws.Cells.Item($r, 7).Value = $row[2]
  $ws.Cells.Item($r, 8).Value = $row[3]
  $ws.Cells.Item($r, 9).Value = $row[4]
}

$ws.Range("L16").Select()
